$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 3.1
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("S7").Value = 1.73
$ws.Range("T7").Value = 2.08
$ws.Range("X7").Value = 8.5
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 26
$ws.Range("AH7").Value = 17
$ws.Range("AL7").Value = 67
$ws.Range("AO7").Value = 15
$ws.Range("AT7").Value = 2
$ws.Range("G12").Value = 1.7
$ws.Range("H12").Value = 3.75
$ws.Range("I12").Value = 4.25
$ws.Range("J12").Value = 2.22
$ws.Range("K12").Value = 2.25
$ws.Range("L12").Value = 4.45
$ws.Range("M12").Value = 1.02
$ws.Range("N12").Value = 12
$ws.Range("R12").Value = 2.05
$ws.Range("S12").Value = 1.3
$ws.Range("T12").Value = 3.38
$ws.Range("U12").Value = 1.6
$ws.Range("V12").Value = 2.07
$ws.Range("W12").Value = 8.5
$ws.Range("X12").Value = 9.25
$ws.Range("Y12").Value = 8
$ws.Range("AA12").Value = 12.5
$ws.Range("AB12").Value = 21
$ws.Range("AE12").Value = 13.5
$ws.Range("AF12").Value = 50
$ws.Range("AG12").Value = 14
$ws.Range("AI12").Value = 14
$ws.Range("AJ12").Value = 70
$ws.Range("AK12").Value = 37
$ws.Range("AL12").Value = 37
$ws.Range("AT12").Value = 3
$ws.Range("AW12").Value = 6.1
$ws.Range("AX12").Value = 23
$ws.Range("AY12").Value = 26
$ws.Range("AZ12").Value = 120
$ws.Range("BA12").Value = 150
$ws.Range("G34").Value = 2.42
$ws.Range("H34").Value = 3.15
$ws.Range("I34").Value = 2.7
$ws.Range("J34").Value = 3
$ws.Range("L34").Value = 3.3
$ws.Range("N34").Value = 6.9
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = 1.75
$ws.Range("X34").Value = 12
$ws.Range("Z34").Value = 26
$ws.Range("AA34").Value = 20
$ws.Range("AB34").Value = 30
$ws.Range("AC34").Value = 6.9
$ws.Range("AD34").Value = 6.2
$ws.Range("AH34").Value = 13.5
$ws.Range("AI34").Value = 10
$ws.Range("AK34").Value = 24
$ws.Range("AL34").Value = 32
$ws.Range("AN34").Value = 4.4
$ws.Range("AO34").Value = 13
$ws.Range("AP34").Value = 20
$ws.Range("AQ34").Value = 55
$ws.Range("AU34").Value = 6.9
$ws.Range("AW34").Value = 4.7
$ws.Range("AX34").Value = 15
$ws.Range("AY34").Value = 22
$ws.Range("AZ34").Value = 65
$ws.Range("BA34").Value = 100
